$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header in H1, copying the formatting from the existing
# header cell G1 so it keeps the bold/bordered/centered header style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" column for the data rows.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 1
